# Insert a new weekly price record for "Ajo" (garlic) at row 154.
# This pushes the existing rows 154-179 down to 155-180 and adds a brand
# new row 154 with the latest week's data (date 2022-06-10 / serial 44722).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 154; Excel copies the formatting
# (including the date style used in column D) from the row above by default.
$ws.Rows.Item(154).Insert()

$ws.Range("A154").Value = 11
$ws.Range("B154").Value = "Vega Monumental Concepción"
$ws.Range("C154").Value = "Bíobío"
$ws.Range("D154").Value = 44722
$ws.Range("E154").Value = 8
$ws.Range("F154").Value = 100112003
$ws.Range("G154").Value = "Ajo"
$ws.Range("H154").Value = "Chino"
$ws.Range("I154").Value = "Primera"
$ws.Range("J154").Value = 220
$ws.Range("K154").Value = 16000
$ws.Range("L154").Value = 17000
$ws.Range("M154").Value = 16545
$ws.Range("N154").Value = "$/caja 10 kilos"
$ws.Range("O154").Value = "China"
$ws.Range("P154").Value = 1654
$ws.Range("Q154").Value = 10
$ws.Range("R154").Value = "Hortaliza"
